$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with uniform run formatting across all runs) ---
$ws.Range("A8").Value = "Volume 30   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/10/2023  Through  4/16/2023"

# --- Crime statistics table updates (Week to Date / 28 Day / YTD / comparison columns) ---
$ws.Range("N14").Copy($ws.Range("M14"))
$ws.Range("M14").Value = -100
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -6.666666666666
$ws.Range("I16").Value = 40
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = -11.111111111111
$ws.Range("L16").Value = 53.846153846153
$ws.Range("M16").Value = 81.818181818181
$ws.Range("N16").Value = -78.021978021978
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 55.555555555555
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = 5.882352941176
$ws.Range("L17").Value = 33.333333333333
$ws.Range("M17").Value = -14.285714285714
$ws.Range("N17").Value = -48.571428571428
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 500
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -15.384615384615
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = -16.666666666666
$ws.Range("L18").Value = -18.918918918918
$ws.Range("M18").Value = -6.25
$ws.Range("N18").Value = -77.611940298507
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 128.571428571429
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 24.390243902439
$ws.Range("I19").Value = 203
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 12.777777777777
$ws.Range("L19").Value = 73.504273504273
$ws.Range("M19").Value = 17.341040462427
$ws.Range("N19").Value = -16.115702479338
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 400
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -5.882352941176
$ws.Range("L20").Value = 300
$ws.Range("M20").Value = 6.666666666666
$ws.Range("N20").Value = -89.873417721519
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 68.421052631578
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = 25.925925925925
$ws.Range("I21").Value = 327
$ws.Range("J21").Value = 317
$ws.Range("K21").Value = 3.154574132492
$ws.Range("L21").Value = 53.521126760563
$ws.Range("M21").Value = 14.335664335664
$ws.Range("N21").Value = -58.867924528301
$ws.Range("J14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E26").Copy($ws.Range("E22"))
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 12
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 300
$ws.Range("C23").Value = 2
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E26").Copy($ws.Range("E23"))
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 14
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = -36.363636363636
$ws.Range("M23").Value = -6.666666666666
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 60
$ws.Range("F24").Value = 56
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = 21.739130434782
$ws.Range("I24").Value = 193
$ws.Range("J24").Value = 180
$ws.Range("K24").Value = 7.222222222222
$ws.Range("L24").Value = 41.911764705882
$ws.Range("M24").Value = -13.839285714285
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -25
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = -16.666666666666
$ws.Range("I25").Value = 101
$ws.Range("J25").Value = 74
$ws.Range("K25").Value = 36.486486486486
$ws.Range("L25").Value = 53.030303030303
$ws.Range("M25").Value = 12.222222222222
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 19
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 26.666666666666
$ws.Range("L27").Value = 171.428571428571
$ws.Range("J14").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 1
$ws.Range("J14").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("N14").Copy($ws.Range("E30"))
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -33.333333333333
$ws.Range("I30").Value = 8
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = 33.333333333333
$ws.Range("N14").Copy($ws.Range("L30"))
$ws.Range("L30").Value = 700
